$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub3 = [char]0x2083

$ws.Range("D2").Value = "57.853.84"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "2.458.74"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'510.89"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").Value = "'133.91"
$ws.Range("E6").Value = "  +3.00%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.557"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").Value = "2.458.52"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "'4.64"
$ws.Range("E13").Value = "  -6.61%  "
$ws.Range("D14").Value = "2.893.53"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "57.920.95"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "'21.98"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("E17").Value = "  +1.68%  "
$ws.Range("D18").Value = "2.414.71"
$ws.Range("E18").Value = "  -1.51%  "
$ws.Range("D19").Value = "'10.37"
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").Value = "'4.17"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "'315.43"
$ws.Range("E21").Value = "  +0.77%  "
$ws.Range("D22").Value = "'6.46"
$ws.Range("E22").Value = "  +5.05%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'5.74"
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("D25").Value = "'65.34"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  -0.58%  "
$ws.Range("D28").Value = "'0.381"
$ws.Range("E28").Value = "  -5.39%  "
$ws.Range("E29").Value = "  +4.57%  "
$ws.Range("D30").Value = "'172.61"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").Value = "0.0{0}0735" -f $sub3
$ws.Range("E31").Value = "  +0.01%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "'1.14"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "'18.11"
$ws.Range("E37").Value = "  +1.61%  "
$ws.Range("D38").Value = "'1.24"
$ws.Range("E38").Value = "  +5.24%  "
$ws.Range("D39").Value = "'3.88"
$ws.Range("E39").Value = "  +2.49%  "
$ws.Range("D40").Value = "'36.78"
$ws.Range("E40").Value = "  +1.30%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("D42").Value = "'0.810"
$ws.Range("E42").Value = "  +0.15%  "
$ws.Range("D43").Value = "'137.20"
$ws.Range("E43").Value = "  +9.40%  "
$ws.Range("D44").Value = "'3.42"
$ws.Range("E44").Value = "  +1.09%  "
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("E46").Value = "  -1.04%  "
$ws.Range("D47").Value = "'256.77"
$ws.Range("E47").Value = "  -0.78%  "
$ws.Range("D48").Value = "'0.0921"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +1.98%  "
$ws.Range("D51").Value = "'17.23"
$ws.Range("E51").Value = "  +1.21%  "
